# TC23_CDS_Filter_PHSAccession-phs002517_FileType_LibStrat_LibSrc.xlsx
#
# The underlying "ParticipantsTab" / "FilesTab" query text and layout are
# unchanged by this revision (the workbook was simply re-opened/re-saved
# by a newer build of Excel, which is why the raw OOXML shows a lot of
# cosmetic churn - shared-string bookkeeping, view metrics, rounding of
# column widths, etc. - none of which alters any cell's visible content).
#
# The one deliberate, user-visible action captured by this commit is that
# the author left the workbook with a different cell selected/active:
# previously C4 was selected (scrolled so column B was the left-most
# visible column, anchored near row 4); afterwards D3 is selected, with
# the view scrolled back up to the top of the sheet.
#
# Reproduce that by activating the sheet, scrolling the window back to
# the top-left area, and selecting D3 - exactly what a user clicking
# around the workbook before saving would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Make sure we are viewing the sheet from the top again (matches
# topLeftCell moving from B4 back to B1) before placing the selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Keep the zoom level the same (70%) - unchanged by this revision.
$excel.ActiveWindow.Zoom = 70

# Move the selection/active cell from C4 to D3.
$ws.Range("D3").Select()
